$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column D entirely (header FACTORS_3 + its 11 data values)
$ws.Range("D1:D12").Delete()

# Update the FACTORS_1 (B) and FACTORS_2 (C) values for the new fleet
$ws.Range("B2").Value = 0.07420953862125235
$ws.Range("C2").Value = 0.05024566007685418

$ws.Range("B3").Value = 0.02681802474505173
$ws.Range("C3").Value = 0.0805501369649695

$ws.Range("B4").Value = 0.01095802025131672
$ws.Range("C4").Value = 0.01392820978207283

$ws.Range("B5").Value = 0.09491908165544634
$ws.Range("C5").Value = 0.1052733130998613

$ws.Range("B6").Value = 0.1728148818351061
$ws.Range("C6").Value = 0.04317076313937656

$ws.Range("B7").Value = 0.140133896904935
$ws.Range("C7").Value = 0.1388960739136928

$ws.Range("B8").Value = 0.1037968923310682
$ws.Range("C8").Value = 0.02661676147595037

$ws.Range("B9").Value = 0.1259133810640778
$ws.Range("C9").Value = 0.1739457295677783

$ws.Range("B10").Value = 0.1509661593483238
$ws.Range("C10").Value = 0.08775544229858545

$ws.Range("B11").Value = 0.0463592010248985
$ws.Range("C11").Value = 0.1537747724949112

$ws.Range("B12").Value = 0.05311092221852348
$ws.Range("C12").Value = 0.1258431371859476
